# Updates the cryptos price/volume table to the latest scraped values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Row = 2;  D = "61.901.06";  E = "  -0.83%  " },
    @{ Row = 3;  D = "2.399.73";   E = "  -1.20%  " },
    @{ Row = 4;  D = $null;        E = "  -0.05%  " },
    @{ Row = 5;  D = "561.70";     E = "  +0.84%  " },
    @{ Row = 6;  D = "142.26";     E = "  -1.27%  " },
    @{ Row = 7;  D = $null;        E = "  +0.08%  " },
    @{ Row = 8;  D = $null;        E = "  -1.04%  " },
    @{ Row = 9;  D = "0.108";      E = "  -1.42%  " },
    @{ Row = 10; D = $null;        E = "  -1.90%  " },
    @{ Row = 11; D = "5.25";       E = "  -2.76%  " },
    @{ Row = 12; D = "0.348";      E = "  -1.32%  " },
    @{ Row = 13; D = "25.51";      E = "  -3.45%  " },
    @{ Row = 14; D = $null;        E = "  -2.07%  " },
    @{ Row = 15; D = "2.832.70";   E = "  -1.29%  " },
    @{ Row = 16; D = "61.900.85";  E = "  -0.30%  " },
    @{ Row = 17; D = "2.381.91";   E = "  -1.88%  " },
    @{ Row = 18; D = "11.20";      E = "  +0.62%  " },
    @{ Row = 19; D = "320.66";     E = "  -1.39%  " },
    @{ Row = 20; D = $null;        E = "  -1.27%  " },
    @{ Row = 21; D = "6.78";       E = "  +0.02%  " },
    @{ Row = 22; D = "0.999";      E = "  -0.18%  " },
    @{ Row = 23; D = "65.93";      E = "  +1.19%  " },
    @{ Row = 24; D = $null;        E = "  -2.02%  " },
    @{ Row = 25; D = "8.78";       E = "  -4.51%  " },
    @{ Row = 26; D = "560.47";     E = "  -2.63%  " },
    @{ Row = 27; D = $null;        E = "  -0.23%  " },
    @{ Row = 28; D = "2.520.12";   E = "  -0.72%  " },
    @{ Row = 29; D = "0.0₃0930";   E = "  -1.90%  " },
    @{ Row = 30; D = "8.13";       E = "  -3.27%  " },
    @{ Row = 31; D = $null;        E = "  -5.02%  " },
    @{ Row = 32; D = "0.146";      E = "  -1.87%  " },
    @{ Row = 33; D = $null;        E = "  -0.48%  " },
    @{ Row = 34; D = "1.49";       E = "  -5.07%  " },
    @{ Row = 35; D = "0.999";      E = $null },
    @{ Row = 36; D = "4.64";       E = "  -3.88%  " },
    @{ Row = 37; D = "151.87";     E = $null },
    @{ Row = 38; D = "5.40";       E = "  -6.24%  " },
    @{ Row = 39; D = $null;        E = "  -2.39%  " },
    @{ Row = 40; D = "18.53";      E = "  -1.85%  " },
    @{ Row = 41; D = "1.78";       E = "  -7.13%  " },
    @{ Row = 42; D = $null;        E = "  -0.02%  " },
    @{ Row = 43; D = "2.24";       E = "  -3.74%  " },
    @{ Row = 44; D = "147.09";     E = "  -3.36%  " },
    @{ Row = 45; D = "3.59";       E = "  -1.58%  " },
    @{ Row = 46; D = $null;        E = "  -3.21%  " },
    @{ Row = 47; D = "19.72";      E = "  -3.95%  " },
    @{ Row = 48; D = "0.585";      E = "  -1.14%  " },
    @{ Row = 49; D = "0.0916";     E = "  +0.10%  " },
    @{ Row = 50; D = $null;        E = "  -2.20%  " },
    @{ Row = 51; D = "11.52";      E = "  +0.32%  " }
)

foreach ($u in $updates) {
    if ($null -ne $u.D) {
        # Prefix with an apostrophe so Excel stores these as text (matching
        # the source data, which keeps prices as literal strings, e.g.
        # "561.70" rather than the number 561.7).
        $ws.Range("D" + $u.Row).Value = "'" + $u.D
    }
    if ($null -ne $u.E) {
        $ws.Range("E" + $u.Row).Value = $u.E
    }
}
